$wb = $excel.ActiveWorkbook

# Update the user's name on the "Users" sheet
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "Drew Koecher"

# Select the sheet and move the active cell selection to C10
$usersSheet.Activate()
$usersSheet.Range("C10").Select()
